# Remove the small duplicate-title "TextBox N" footer shapes that were
# left over on every slide (TITLE, SYSTEM FLOW, REFERENCES, ABSTRACT,
# INTRODUCTION, LITRATURE SURVEY, EXISTING SYSTEM, WORKING, ...).
$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -like "TextBox *") {
            $sh.Delete()
        }
    }
}
